# Edit script: Add "Socks in the Dark" problem (intro + part A) after the
# "A Cat, a Parrot, and a Bag of Seed" problem, fix up 3 header1.xml sdt
# blocks (add empty <w:sdtEndPr/>), and merge split date/course runs in
# header2.xml into a single run.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Body content: insert the "Socks in the Dark" intro + two list items
#    and the "Define." paragraph right after the existing
#    "This solution will work..." paragraph, reusing the final (empty)
#    trailing paragraph for the new "Define." text and relocating the
#    _GoBack bookmark there too.
# ---------------------------------------------------------------------
$paras = $d.Paragraphs
$pLastSolution = $paras[17]
$pTrailingEmpty = $paras[18]
$bodyRange = $d.Range($pLastSolution.Range.Start, $pTrailingEmpty.Range.End)
$bodyRange.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="009E0211" w:rsidRPr="009E0211" w:rsidRDefault="009E0211" w:rsidP="009E0211"> <w:pPr> <w:pStyle w:val="NormalWeb"/> <w:rPr> <w:sz w:val="22"/> <w:szCs w:val="22"/> </w:rPr> </w:pPr> <w:r w:rsidRPr="009E0211"> <w:rPr> <w:sz w:val="22"/> <w:szCs w:val="22"/> </w:rPr> <w:t>This solution will work since there are no rules or constraints against switching items/creatures, complying with the rule to transporting one creature/item one at a time on the boat.</w:t> </w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"> <w:pPr> <w:pStyle w:val="NormalWeb"/> </w:pPr> <w:r> <w:rPr> <w:rFonts w:ascii="Hero" w:hAnsi="Hero"/> <w:color w:val="7F7F7F"/> <w:sz w:val="24"/> <w:szCs w:val="24"/> </w:rPr> <w:t xml:space="preserve">Socks in the Dark: </w:t> </w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"> <w:pPr> <w:pStyle w:val="NormalWeb"/> </w:pPr> <w:r> <w:rPr> <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/> <w:sz w:val="24"/> <w:szCs w:val="24"/> </w:rPr> <w:t xml:space="preserve">There are 20 socks in a drawer: 5 pairs of black socks, 3 pairs of brown and 2 pairs of white. You select the socks in the dark and can check them only after a selection has been made. What is the smallest number of socks you need to select to guarantee getting the </w:t> </w:r> <w:proofErr w:type="gramStart"/> <w:r> <w:rPr> <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/> <w:sz w:val="24"/> <w:szCs w:val="24"/> </w:rPr> <w:t>following:</w:t> </w:r> <w:proofErr w:type="gramEnd"/> <w:r> <w:rPr> <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/> <w:sz w:val="24"/> <w:szCs w:val="24"/> </w:rPr> <w:t xml:space="preserve"> </w:t> </w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"> <w:pPr> <w:pStyle w:val="NormalWeb"/> <w:numPr> <w:ilvl w:val="0"/> <w:numId w:val="10"/> </w:numPr> </w:pPr> <w:r> <w:rPr> <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/> <w:sz w:val="24"/> <w:szCs w:val="24"/> </w:rPr> <w:lastRenderedPageBreak/> <w:t>At least one matching pair</w:t> </w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"> <w:pPr> <w:pStyle w:val="NormalWeb"/> <w:numPr> <w:ilvl w:val="0"/> <w:numId w:val="10"/> </w:numPr> </w:pPr> <w:r> <w:rPr> <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/> <w:sz w:val="24"/> <w:szCs w:val="24"/> </w:rPr> <w:t xml:space="preserve">At least one matching pair </w:t> </w:r> <w:r> <w:rPr> <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/> <w:i/> <w:iCs/> <w:sz w:val="24"/> <w:szCs w:val="24"/> </w:rPr> <w:t>of each color.</w:t> </w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="003027E2" w:rsidRPr="003B334B" w:rsidRDefault="003027E2" w:rsidP="003027E2"> <w:pPr> <w:pStyle w:val="NormalWeb"/> <w:rPr> <w:sz w:val="22"/> <w:szCs w:val="22"/> </w:rPr> </w:pPr> <w:r> <w:rPr> <w:b/> </w:rPr> <w:t xml:space="preserve">Define. </w:t> </w:r> <w:r> <w:t>The first part of the problem is to find a matching pair of socks, regardless of color.</w:t> </w:r> <w:bookmarkStart w:id="0" w:name="_GoBack"/> <w:bookmarkEnd w:id="0"/></w:p>')

# ---------------------------------------------------------------------
# 2) header1.xml (the unused "even pages" header): add an empty
#    <w:sdtEndPr/> to each of the 3 placeholder content controls.
# ---------------------------------------------------------------------
$sec = $d.Sections[1]
$evenHeader = $sec.Headers[3]
$evenHeader.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w14:paraId="11F31F79" w14:textId="77777777" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w:rsidR="00405106" w:rsidRDefault="00405106"><w:pPr><w:pStyle w:val="Header"/></w:pPr><w:sdt><w:sdtPr><w:id w:val="171999623"/><w:placeholder><w:docPart w:val="D4EE7E02F168FA43A45A0A2731881CA1"/></w:placeholder><w:temporary/><w:showingPlcHdr/></w:sdtPr><w:sdtEndPr/><w:sdtContent><w:r><w:t>[Type text]</w:t></w:r></w:sdtContent></w:sdt><w:r><w:ptab w:relativeTo="margin" w:alignment="center" w:leader="none"/></w:r><w:sdt><w:sdtPr><w:id w:val="171999624"/><w:placeholder><w:docPart w:val="640AB2700C8AFA439924544B003FD028"/></w:placeholder><w:temporary/><w:showingPlcHdr/></w:sdtPr><w:sdtEndPr/><w:sdtContent><w:r><w:t>[Type text]</w:t></w:r></w:sdtContent></w:sdt><w:r><w:ptab w:relativeTo="margin" w:alignment="right" w:leader="none"/></w:r><w:sdt><w:sdtPr><w:id w:val="171999625"/><w:placeholder><w:docPart w:val="505DB3F4FA421B4F81E5679EB084935C"/></w:placeholder><w:temporary/><w:showingPlcHdr/></w:sdtPr><w:sdtEndPr/><w:sdtContent><w:r><w:t>[Type text]</w:t></w:r></w:sdtContent></w:sdt></w:p>')

# ---------------------------------------------------------------------
# 3) header2.xml (the default/primary header): merge the "7/9/14",
#    the run of spaces, and the "Web Programming..." run into one run.
# ---------------------------------------------------------------------
$primaryHeader = $sec.Headers[1]
$primaryHeader.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="701D062C" w14:textId="1B4AE59A" w:rsidR="00405106" w:rsidRDefault="00405106"><w:pPr><w:pStyle w:val="Header"/></w:pPr><w:r w:rsidRPr="009E0C0A"><w:rPr><w:b/></w:rPr><w:t>Cathleen Carbonell</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve">7/9/14       Web Programming Fundamentals Section 1 </w:t></w:r></w:p>')

Write-Host "Edit complete."
